$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLOCK")
$line = $wb.Worksheets.Add($ws)
$line.Name = "LINE"
$section = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $line)
$section.Name = "SECTION"

# LINE sheet content
$line.Range("A1").Value = "line"
$line.Range("A2").Value = "GREEN"
$line.Range("A3").Value = "RED"

$headerFormatSrc = $wb.Worksheets.Item("STATION").Range("A1")
$headerFormatSrc.Copy()
$line.Range("A1").PasteSpecial(-4122)
$line.PageSetup.Orientation = 1

# SECTION sheet content
$section.Range("A1").Value = "line"
$section.Range("B1").Value = "section"

$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","ZZ","YY")
$row = 2
for ($i = 0; $i -lt 28; $i++) {
    $section.Cells.Item($row, 1).Value = "GREEN"
    $section.Cells.Item($row, 2).Value = $letters[$i]
    $row++
}
for ($i = 0; $i -lt 21; $i++) {
    $section.Cells.Item($row, 1).Value = "RED"
    $section.Cells.Item($row, 2).Value = $letters[$i]
    $row++
}

$headerFormatSrc.Copy()
$section.Range("A1:B1").PasteSpecial(-4122)

# BLOCK sheet fixes
$block = $wb.Worksheets.Item("BLOCK")
$block.Range("C5").Copy()
$block.Range("E5").PasteSpecial(-4122)
[void]$block.Range("L10").Select()

# Restore view/selection state: LINE's selection sits on A4 (next empty
# row), SECTION's selection sits on D5, and SECTION remains the active
# (frontmost) sheet, matching the workbook's saved activeTab.
[void]$line.Range("A4").Select()
[void]$section.Range("D5").Select()
